$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextValue "D2" "26.538.61"
Set-TextValue "E2" "  -2.75%  "
Set-TextValue "D3" "1.668.63"
Set-TextValue "E3" "  -2.23%  "
Set-TextValue "E4" "  +0.56%  "
Set-TextValue "D5" "219.12"
Set-TextValue "E5" "  -1.82%  "
Set-TextValue "D6" "0.5123"
Set-TextValue "E6" "  -3.39%  "
Set-TextValue "E7" "  +0.54%  "
Set-TextValue "D8" "0.06431"
Set-TextValue "E8" "  -2.31%  "
Set-TextValue "D9" "0.2557"
Set-TextValue "E9" "  -3.68%  "
Set-TextValue "D10" "19.89"
Set-TextValue "E10" "  -4.68%  "
Set-TextValue "D11" "0.07640"
Set-TextValue "E11" "  -0.05%  "
Set-TextValue "D12" "4.333"
Set-TextValue "E12" "  -5.53%  "
Set-TextValue "B13" "WrappedliquidstakedEther2.0"
Set-TextValue "C13" "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
Set-TextValue "D13" "1.897.64"
Set-TextValue "E13" "  -2.25%  "
Set-TextValue "B14" "WrappedEther"
Set-TextValue "C14" "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue "D14" "1.669.18"
Set-TextValue "E14" "  -2.54%  "
Set-TextValue "D15" "0.5551"
Set-TextValue "E15" "  -3.23%  "
Set-TextValue "D16" "0.0₅7990"
Set-TextValue "E16" "  -2.38%  "
Set-TextValue "D17" "64.44"
Set-TextValue "E17" "  -4.50%  "
Set-TextValue "D18" "26.568.92"
Set-TextValue "E18" "  -2.56%  "
Set-TextValue "E19" "  +0.53%  "
Set-TextValue "D20" "208.25"
Set-TextValue "E20" "  -3.75%  "
Set-TextValue "D21" "4.439"
Set-TextValue "E21" "  -5.10%  "
Set-TextValue "E22" "  -3.63%  "
Set-TextValue "D23" "5.871"
Set-TextValue "E23" "  -1.73%  "
Set-TextValue "D24" "1.009"
Set-TextValue "E24" "  +0.65%  "
Set-TextValue "D25" "142.89"
Set-TextValue "E25" "  +0.47%  "
Set-TextValue "D26" "1.723"
Set-TextValue "E26" "  -1.17%  "
Set-TextValue "D28" "6.968"
Set-TextValue "E28" "  -3.96%  "
Set-TextValue "D29" "15.63"
Set-TextValue "E29" "  -4.19%  "
Set-TextValue "D30" "0.05200"
Set-TextValue "E30" "  -3.18%  "
Set-TextValue "D31" "1.263"
Set-TextValue "E31" "  -2.12%  "
Set-TextValue "D32" "3.339"
Set-TextValue "E32" "  -4.82%  "
Set-TextValue "D33" "3.180"
Set-TextValue "E33" "  -7.09%  "
Set-TextValue "D34" "1.577"
Set-TextValue "E34" "  -3.66%  "
Set-TextValue "D35" "2.759"
Set-TextValue "E35" "  -4.06%  "
Set-TextValue "D36" "2.375"
Set-TextValue "E36" "  -1.80%  "
Set-TextValue "D37" "0.9215"
Set-TextValue "E37" "  -2.66%  "
Set-TextValue "D38" "0.5735"
Set-TextValue "E38" "  -2.50%  "
Set-TextValue "D39" "1.156.71"
Set-TextValue "E39" "  +11.25%  "
Set-TextValue "D40" "0.01586"
Set-TextValue "E40" "  -2.95%  "
Set-TextValue "E41" "  +0.54%  "
Set-TextValue "B42" "TrustWalletToken"
Set-TextValue "C42" "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextValue "D42" "0.8295"
Set-TextValue "E42" "  -1.24%  "
Set-TextValue "B43" "FraxShare"
Set-TextValue "C43" "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue "D43" "5.650"
Set-TextValue "E43" "  -3.54%  "
Set-TextValue "D44" "99.89"
Set-TextValue "E44" "  -1.13%  "
Set-TextValue "D45" "1.807.35"
Set-TextValue "E45" "  -2.21%  "
Set-TextValue "D46" "0.0₈111"
Set-TextValue "E46" "  -3.84%  "
Set-TextValue "D47" "0.4491"
Set-TextValue "E47" "  -0.14%  "
Set-TextValue "D48" "55.49"
Set-TextValue "E48" "  -4.51%  "
Set-TextValue "E49" "  -0.20%  "
Set-TextValue "D50" "7.933"
Set-TextValue "E50" "  -1.89%  "
Set-TextValue "D51" "0.05140"
Set-TextValue "E51" "  -1.91%  "
